# Delete the "Type" column (column C) entirely from the CapitalCommitment sheet.
# This shifts all subsequent columns (D.. ) one position to the left, and Excel
# automatically adjusts the dimension, data validations, comments, hyperlinks
# and column widths that referenced those columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Columns.Item(3).Select()
$ws.Columns.Item(3).Delete()
